# Add file upload functionality: append 3 new daily records (rows 124-126)
# to each of the 4 data sheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2).

$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Index = 1
        Rows = @(
            @{ A = 45910.46265046296; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x00,0xF0"; E = "0x07"; F = 400;  G = [double]"5.68631262647113e+23"; H = 240; I = 7 },
            @{ A = 45911.46655092593; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x00,0xEC"; E = "0x07"; F = 400;  G = [double]"5.68631262647113e+23"; H = 240; I = 7 },
            @{ A = 45912.46122685185; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x00,0xEC"; E = "0x07"; F = 400;  G = [double]"5.68631262647113e+23"; H = 236; I = 7 }
        )
    },
    @{
        Index = 2
        Rows = @(
            @{ A = 45910.46265046296; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x04"; E = "0x19"; F = 380;  G = [double]"5.68432987514711e+23"; H = 260; I = 25 },
            @{ A = 45911.46655092593; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x00"; E = "0x19"; F = 380;  G = [double]"5.68432987514711e+23"; H = 256; I = 25 },
            @{ A = 45912.46122685185; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x00"; E = "0x19"; F = 380;  G = [double]"5.68432987514711e+23"; H = 256; I = 25 }
        )
    },
    @{
        Index = 3
        Rows = @(
            @{ A = 45910.46265046296; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x58"; E = "0x15"; F = 110;  G = [double]"5.68631262647113e+23"; H = 88;  I = 15 },
            @{ A = 45911.46655092593; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x58"; E = "0x15"; F = 110;  G = [double]"5.68631262647113e+23"; H = 88;  I = 15 },
            @{ A = 45912.46122685185; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x58"; E = "0x15"; F = 110;  G = [double]"5.68631262647113e+23"; H = 88;  I = 15 }
        )
    },
    @{
        Index = 4
        Rows = @(
            @{ A = 45910.46265046296; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x6F"; E = "0x9"; F = 130;  G = [double]"5.68631262647113e+23"; H = 111; I = 9 },
            @{ A = 45911.46655092593; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x6F"; E = "0x9"; F = 130;  G = [double]"5.68631262647113e+23"; H = 111; I = 9 },
            @{ A = 45912.46122685185; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x6E"; E = "0x9"; F = 130;  G = [double]"5.68631262647113e+23"; H = 110; I = 9 }
        )
    }
)

foreach ($sheetInfo in $sheetsData) {
    $ws = $wb.Worksheets.Item($sheetInfo.Index)
    $dateFormat = $ws.Cells.Item(123, 1).NumberFormat
    $rowNum = 124
    foreach ($r in $sheetInfo.Rows) {
        $ws.Cells.Item($rowNum, 1).Value = $r.A
        $ws.Cells.Item($rowNum, 1).NumberFormat = $dateFormat
        $ws.Cells.Item($rowNum, 2).Value = $r.B
        $ws.Cells.Item($rowNum, 3).Value = $r.C
        $ws.Cells.Item($rowNum, 4).Value = $r.D
        $ws.Cells.Item($rowNum, 5).Value = $r.E
        $ws.Cells.Item($rowNum, 6).Value = $r.F
        $ws.Cells.Item($rowNum, 7).Value = $r.G
        $ws.Cells.Item($rowNum, 8).Value = $r.H
        $ws.Cells.Item($rowNum, 9).Value = $r.I
        $rowNum = $rowNum + 1
    }
}
